$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.14"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.83%"
$ws.Range("E2").NumberFormat = "General"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.12"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.04%"
$ws.Range("E3").NumberFormat = "General"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.154"
$ws.Range("D4").NumberFormat = "General"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08238"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "5.09%"
$ws.Range("E5").NumberFormat = "General"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.152"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.03%"
$ws.Range("E6").NumberFormat = "General"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.012"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.97%"
$ws.Range("E7").NumberFormat = "General"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.135"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.76%"
$ws.Range("E8").NumberFormat = "General"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9272"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.09%"
$ws.Range("E9").NumberFormat = "General"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1007"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.23%"
$ws.Range("E10").NumberFormat = "General"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1892"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.63%"
$ws.Range("E11").NumberFormat = "General"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09237"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.87%"
$ws.Range("E12").NumberFormat = "General"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03622"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.52%"
$ws.Range("E13").NumberFormat = "General"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09916"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.18%"
$ws.Range("E14").NumberFormat = "General"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001439"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.30%"
$ws.Range("E15").NumberFormat = "General"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005782"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.70%"
$ws.Range("E16").NumberFormat = "General"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.459"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.12%"
$ws.Range("E17").NumberFormat = "General"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.797"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "13.76%"
$ws.Range("E18").NumberFormat = "General"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.43%"
$ws.Range("E19").NumberFormat = "General"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.211"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "8.08%"
$ws.Range("E20").NumberFormat = "General"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.68%"
$ws.Range("E22").NumberFormat = "General"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04597"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.34%"
$ws.Range("E23").NumberFormat = "General"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.69%"
$ws.Range("E24").NumberFormat = "General"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004733"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.97%"
$ws.Range("E25").NumberFormat = "General"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-21.96%"
$ws.Range("E26").NumberFormat = "General"

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.33%"
$ws.Range("E27").NumberFormat = "General"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02007"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.72%"
$ws.Range("E39").NumberFormat = "General"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04923"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.29%"
$ws.Range("E40").NumberFormat = "General"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007761"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.21%"
$ws.Range("E41").NumberFormat = "General"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1400"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.14%"
$ws.Range("E42").NumberFormat = "General"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007832"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.09%"
$ws.Range("E43").NumberFormat = "General"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002104"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.08%"
$ws.Range("E44").NumberFormat = "General"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01193"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.23%"
$ws.Range("E45").NumberFormat = "General"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006470"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.92%"
$ws.Range("E46").NumberFormat = "General"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("E47").NumberFormat = "General"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.86"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.95%"
$ws.Range("E48").NumberFormat = "General"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.06%"
$ws.Range("E49").NumberFormat = "General"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("E50").NumberFormat = "General"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.08%"
$ws.Range("E51").NumberFormat = "General"
